$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update formulas in row 11 (divide by 4 instead of 5)
$ws.Range("E11").Formula = "=C10-(D10-E10)-((D10-C10)/4)"
$ws.Range("H11").Formula = "=F10-(G10-H10)-((G10-F10)/4)"
$ws.Range("K11").Formula = "=I10-(J10-K10)-((J10-I10)/4)"
$ws.Range("N11").Formula = "=L10-(M10-N10)-((M10-L10)/4)"

# Update formulas in row 12 (divide by 4 instead of 5)
$ws.Range("D12").Formula = "=(D10-C10)/4"
$ws.Range("G12").Formula = "=(G10-F10)/4"
$ws.Range("J12").Formula = "=(J10-I10)/4"
$ws.Range("M12").Formula = "=(M10-L10)/4"

# Update the selection to D12 as in the diff
$ws.Range("D12").Select()
